# Auto-generated Excel COM-interop script
# Applies a scheduled data refresh to the per-job Leve profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all class sheets.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1405.1052
$ws.Range("I70").Value = 889.8
$ws.Range("J70").Value = 1977.6666
$ws.Range("K70").Value = 2669.4
$ws.Range("L70").Value = 5932.9998
$ws.Range("M70").Value = -2399.4
$ws.Range("N70").Value = -6472.9998
$ws.Range("H73").Value = 1405.1052
$ws.Range("I73").Value = 889.8
$ws.Range("J73").Value = 1977.6666
$ws.Range("K73").Value = 2669.4
$ws.Range("L73").Value = 5932.9998
$ws.Range("M73").Value = -1733.4
$ws.Range("N73").Value = -7804.9998
$ws.Range("H132").Value = 1661.3226
$ws.Range("I132").Value = 1611.3462
$ws.Range("J132").Value = 1921.2
$ws.Range("K132").Value = 4834.0386
$ws.Range("L132").Value = 5763.6
$ws.Range("M132").Value = -2304.0386
$ws.Range("N132").Value = -10823.6
$ws.Range("H135").Value = 5651.727
$ws.Range("I135").Value = 5326.6
$ws.Range("J135").Value = 5922.6665
$ws.Range("K135").Value = 47939.4
$ws.Range("L135").Value = 53303.9985
$ws.Range("M135").Value = -45404.4
$ws.Range("N135").Value = -58373.9985
$ws.Range("H137").Value = 13515282
$ws.Range("I137").Value = 1601.6364
$ws.Range("J137").Value = 33335346
$ws.Range("K137").Value = 4804.9092
$ws.Range("L137").Value = 100006038
$ws.Range("M137").Value = -2254.9092
$ws.Range("N137").Value = -100011138
$ws.Range("H138").Value = 4477.66
$ws.Range("I138").Value = 1951.8889
$ws.Range("J138").Value = 4727.4614
$ws.Range("K138").Value = 5855.6667
$ws.Range("L138").Value = 14182.3842
$ws.Range("M138").Value = -715.6666999999998
$ws.Range("N138").Value = -24462.3842
$ws.Range("H139").Value = 45800
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 45800
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 45800
$ws.Range("N139").Value = -56080
$ws.Range("H141").Value = 1834.375
$ws.Range("I141").Value = 1610.7142
$ws.Range("J141").Value = 3400
$ws.Range("K141").Value = 4832.142599999999
$ws.Range("L141").Value = 10200
$ws.Range("M141").Value = 347.8574000000008
$ws.Range("N141").Value = -20560

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24668.227
$ws.Range("I32").Value = 21789.654
$ws.Range("J32").Value = 47285.57
$ws.Range("K32").Value = 21789.654
$ws.Range("L32").Value = 47285.57
$ws.Range("M32").Value = -21502.654
$ws.Range("N32").Value = -47859.57
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H132").Value = 1685.2307
$ws.Range("I132").Value = 1262.0769
$ws.Range("J132").Value = 2531.5386
$ws.Range("K132").Value = 3786.2307
$ws.Range("L132").Value = 7594.6158
$ws.Range("M132").Value = -1256.2307
$ws.Range("N132").Value = -12654.6158
$ws.Range("N39").ClearContents()
$ws.Range("N129").ClearContents()

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("H134").Value = 2397.0476
$ws.Range("I134").Value = 2333.5789
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 7000.736699999999
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -4465.736699999999
$ws.Range("N134").Value = -14070
$ws.Range("M5").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("N59").ClearContents()

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 3505.5
$ws.Range("I12").Value = 2005
$ws.Range("J12").Value = 5006
$ws.Range("K12").Value = 2005
$ws.Range("L12").Value = 5006
$ws.Range("M12").Value = -1835
$ws.Range("N12").Value = -5346
$ws.Range("H39").Value = 3057
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3057
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 3057
$ws.Range("N39").Value = -3839
$ws.Range("H49").Value = 3057
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 3057
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 3057
$ws.Range("N49").Value = -3421
$ws.Range("H132").Value = 2097.3044
$ws.Range("I132").Value = 1404.75
$ws.Range("J132").Value = 3680.2856
$ws.Range("K132").Value = 4214.25
$ws.Range("L132").Value = 11040.8568
$ws.Range("M132").Value = -1684.25
$ws.Range("N132").Value = -16100.8568
$ws.Range("H140").Value = 49763.156
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 49763.156
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 49763.156
$ws.Range("N140").Value = -60123.156
$ws.Range("M39").ClearContents()
$ws.Range("M49").ClearContents()

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("H44").Value = 598.5714
$ws.Range("I44").Value = 580
$ws.Range("J44").Value = 600
$ws.Range("K44").Value = 1740
$ws.Range("L44").Value = 1800
$ws.Range("M44").Value = -1342
$ws.Range("N44").Value = -2596
$ws.Range("H47").Value = 813
$ws.Range("I47").Value = 126.666664
$ws.Range("J47").Value = 1499.3334
$ws.Range("K47").Value = 379.999992
$ws.Range("L47").Value = 4498.0002
$ws.Range("M47").Value = 51.00000799999998
$ws.Range("N47").Value = -5360.0002
$ws.Range("H103").Value = 519.6
$ws.Range("I103").Value = 199.66667
$ws.Range("J103").Value = 999.5
$ws.Range("K103").Value = 599.00001
$ws.Range("L103").Value = 2998.5
$ws.Range("M103").Value = 279.99999
$ws.Range("N103").Value = -4756.5
$ws.Range("H113").Value = 531.24
$ws.Range("I113").Value = 472.72726
$ws.Range("J113").Value = 547.7436
$ws.Range("K113").Value = 1418.18178
$ws.Range("L113").Value = 1643.2308
$ws.Range("M113").Value = 751.8182200000001
$ws.Range("N113").Value = -5983.2308
$ws.Range("H114").Value = 1175.2759
$ws.Range("I114").Value = 481.63635
$ws.Range("J114").Value = 1599.1666
$ws.Range("K114").Value = 1444.90905
$ws.Range("L114").Value = 4797.4998
$ws.Range("M114").Value = 1809.09095
$ws.Range("N114").Value = -11305.4998
$ws.Range("H117").Value = 2255.8667
$ws.Range("I117").Value = 963.2
$ws.Range("J117").Value = 2902.2
$ws.Range("K117").Value = 2889.6
$ws.Range("L117").Value = 8706.599999999999
$ws.Range("M117").Value = 552.3999999999996
$ws.Range("N117").Value = -15590.6
$ws.Range("H118").Value = 5675.8335
$ws.Range("I118").Value = 1226.6666
$ws.Range("J118").Value = 7158.8887
$ws.Range("K118").Value = 3679.9998
$ws.Range("L118").Value = 21476.6661
$ws.Range("M118").Value = -2436.9998
$ws.Range("N118").Value = -23962.6661
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4342.409
$ws.Range("I70").Value = 4201.8335
$ws.Range("J70").Value = 4975
$ws.Range("K70").Value = 4201.8335
$ws.Range("L70").Value = 4975
$ws.Range("M70").Value = -3931.8335
$ws.Range("H73").Value = 4342.409
$ws.Range("I73").Value = 4201.8335
$ws.Range("J73").Value = 4975
$ws.Range("K73").Value = 4201.8335
$ws.Range("L73").Value = 4975
$ws.Range("M73").Value = -3265.8335
$ws.Range("H80").Value = 3380
$ws.Range("I80").Value = 3380
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3380
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2382
$ws.Range("H83").Value = 3380
$ws.Range("I83").Value = 3380
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 16900
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -11908
$ws.Range("H107").Value = 706.4375
$ws.Range("I107").Value = 721.6429000000001
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 721.6429000000001
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 1198.3571
$ws.Range("N107").Value = -4440
$ws.Range("H132").Value = 2098.2144
$ws.Range("I132").Value = 1725.44
$ws.Range("J132").Value = 2646.4119
$ws.Range("K132").Value = 5176.32
$ws.Range("L132").Value = 7939.2357
$ws.Range("M132").Value = -2646.32
$ws.Range("N132").Value = -12999.2357

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 344.6
$ws.Range("I22").Value = 310.5
$ws.Range("J22").Value = 367.33334
$ws.Range("K22").Value = 310.5
$ws.Range("L22").Value = 367.33334
$ws.Range("M22").Value = -15.5
$ws.Range("N22").Value = -957.33334
$ws.Range("H27").Value = 344.6
$ws.Range("I27").Value = 310.5
$ws.Range("J27").Value = 367.33334
$ws.Range("K27").Value = 310.5
$ws.Range("L27").Value = 367.33334
$ws.Range("M27").Value = -203.5
$ws.Range("N27").Value = -581.33334
$ws.Range("H132").Value = 4333.033
$ws.Range("I132").Value = 3465.1765
$ws.Range("J132").Value = 5467.923
$ws.Range("K132").Value = 10395.5295
$ws.Range("L132").Value = 16403.769
$ws.Range("M132").Value = -7865.529500000001
$ws.Range("N132").Value = -21463.769
$ws.Range("H139").Value = 31809.092
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 31809.092
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 31809.092
$ws.Range("N139").Value = -42089.092

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 321.5
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 395.33334
$ws.Range("K23").Value = 100
$ws.Range("L23").Value = 395.33334
$ws.Range("M23").Value = 129
$ws.Range("N23").Value = -853.33334
$ws.Range("H45").Value = 4837.75
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 4837.75
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 4837.75
$ws.Range("N45").Value = -5819.75
$ws.Range("H62").Value = 4501.25
$ws.Range("I62").Value = 4201
$ws.Range("J62").Value = 4801.5
$ws.Range("K62").Value = 4201
$ws.Range("L62").Value = 4801.5
$ws.Range("M62").Value = -3577
$ws.Range("H65").Value = 4501.25
$ws.Range("I65").Value = 4201
$ws.Range("J65").Value = 4801.5
$ws.Range("K65").Value = 21005
$ws.Range("L65").Value = 24007.5
$ws.Range("M65").Value = -17885
